$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a new "Wut" material column to Table6 (E3:J4) on Sheet1.
$ws1.Range("J3").Value = "Wut"
$ws1.Range("J4").Value = -1

# Switch the active tile-material selector from "Stone" to "Glass".
$ws1.Range("H2").Value = "Glass"

# Switch the AtkPoint calc mode from "per second" to "per swing".
$ws1.Range("L2").Value = "AtkPnt per swing"

# Unrelated label swap near B3 (SwingsPerSec note).
$ws1.Range("C3").Value = "Wut"

# Update the selected range shown when the workbook is opened.
$ws1.Activate()
$ws1.Range("A4:C4").Select()
